$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2700
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 2840
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 2840
$ws.Range("M32").Value = -1674
$ws.Range("N32").Value = -3492
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("H137").Value = 2224.8057
$ws.Range("I137").Value = 2965.0625
$ws.Range("J137").Value = 1632.6
$ws.Range("K137").Value = 8895.1875
$ws.Range("L137").Value = 4897.799999999999
$ws.Range("M137").Value = -6345.1875
$ws.Range("N137").Value = -9997.799999999999
$ws.Range("N75").ClearContents()
$ws.Range("N78").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3281.3333
$ws.Range("I45").Value = 1350
$ws.Range("J45").Value = 5037.091
$ws.Range("K45").Value = 1350
$ws.Range("L45").Value = 5037.091
$ws.Range("M45").Value = -973
$ws.Range("N45").Value = -5791.091
$ws.Range("H97").Value = 1008.4231
$ws.Range("I97").Value = 796.9
$ws.Range("K97").Value = 796.9
$ws.Range("M97").Value = -300.9
$ws.Range("H122").Value = 2980
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2980
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 8940
$ws.Range("N122").Value = -13840
$ws.Range("H132").Value = 4609.4653
$ws.Range("I132").Value = 2003.04
$ws.Range("J132").Value = 8229.5
$ws.Range("K132").Value = 6009.12
$ws.Range("L132").Value = 24688.5
$ws.Range("M132").Value = -3479.12
$ws.Range("N132").Value = -29748.5
$ws.Range("M122").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1590
$ws.Range("I20").Value = 580.5714
$ws.Range("J20").Value = 1897.2174
$ws.Range("K20").Value = 580.5714
$ws.Range("L20").Value = 1897.2174
$ws.Range("M20").Value = -333.5714
$ws.Range("N20").Value = -2391.2174
$ws.Range("H30").Value = 3000
$ws.Range("I30").Value = 3000
$ws.Range("K30").Value = 3000
$ws.Range("M30").Value = -2875
$ws.Range("H107").Value = 1620.1786
$ws.Range("I107").Value = 1615.409
$ws.Range("K107").Value = 1615.409
$ws.Range("M107").Value = 304.5909999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 18208.572
$ws.Range("I62").Value = 13321.111
$ws.Range("J62").Value = 27006
$ws.Range("K62").Value = 13321.111
$ws.Range("L62").Value = 27006
$ws.Range("M62").Value = -12697.111
$ws.Range("N62").Value = -28254
$ws.Range("H65").Value = 18208.572
$ws.Range("I65").Value = 13321.111
$ws.Range("J65").Value = 27006
$ws.Range("K65").Value = 66605.55500000001
$ws.Range("L65").Value = 135030
$ws.Range("M65").Value = -63485.55500000001
$ws.Range("N65").Value = -141270
$ws.Range("H107").Value = 1124.6666
$ws.Range("I107").Value = 448.41666
$ws.Range("J107").Value = 2477.1667
$ws.Range("K107").Value = 448.41666
$ws.Range("L107").Value = 2477.1667
$ws.Range("M107").Value = 1471.58334
$ws.Range("N107").Value = -6317.1667
$ws.Range("H122").Value = 100002170
$ws.Range("I122").Value = 166667500
$ws.Range("J122").Value = 4174.75
$ws.Range("K122").Value = 500002500
$ws.Range("L122").Value = 12524.25
$ws.Range("M122").Value = -500000050
$ws.Range("N122").Value = -17424.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 6000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 6000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 18000
$ws.Range("N24").Value = -18460
$ws.Range("H34").Value = 14818.875
$ws.Range("J34").Value = 15793.467
$ws.Range("L34").Value = 47380.401
$ws.Range("N34").Value = -47548.401
$ws.Range("H87").Value = 6144.857
$ws.Range("I87").Value = 1014
$ws.Range("J87").Value = 7000
$ws.Range("K87").Value = 3042
$ws.Range("L87").Value = 21000
$ws.Range("M87").Value = -1794
$ws.Range("N87").Value = -23496
$ws.Range("H90").Value = 6144.857
$ws.Range("I90").Value = 1014
$ws.Range("J90").Value = 7000
$ws.Range("K90").Value = 9126
$ws.Range("L90").Value = 63000
$ws.Range("M90").Value = -2886
$ws.Range("N90").Value = -75480
$ws.Range("H113").Value = 761.37036
$ws.Range("I113").Value = 595.5833
$ws.Range("J113").Value = 894
$ws.Range("K113").Value = 1786.7499
$ws.Range("L113").Value = 2682
$ws.Range("M113").Value = 383.2501
$ws.Range("N113").Value = -7022
$ws.Range("H129").Value = 2062.5293
$ws.Range("I129").Value = 883.3333
$ws.Range("J129").Value = 2315.2144
$ws.Range("K129").Value = 2649.9999
$ws.Range("L129").Value = 6945.6432
$ws.Range("M129").Value = 2350.0001
$ws.Range("N129").Value = -16945.6432
$ws.Range("H132").Value = 2691.3635
$ws.Range("J132").Value = 3089.4443
$ws.Range("L132").Value = 27804.9987
$ws.Range("N132").Value = -32864.9987
$ws.Range("M24").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 10000500
$ws.Range("I29").Value = 10000500
$ws.Range("K29").Value = 10000500
$ws.Range("M29").Value = -10000210
$ws.Range("H70").Value = 5488.278
$ws.Range("I70").Value = 4990.7144
$ws.Range("J70").Value = 5804.909
$ws.Range("K70").Value = 4990.7144
$ws.Range("L70").Value = 5804.909
$ws.Range("M70").Value = -4720.7144
$ws.Range("N70").Value = -6344.909
$ws.Range("H73").Value = 5488.278
$ws.Range("I73").Value = 4990.7144
$ws.Range("J73").Value = 5804.909
$ws.Range("K73").Value = 4990.7144
$ws.Range("L73").Value = 5804.909
$ws.Range("M73").Value = -4054.7144
$ws.Range("N73").Value = -7676.909
$ws.Range("H132").Value = 2936.5173
$ws.Range("I132").Value = 2923.8572
$ws.Range("J132").Value = 2940.5454
$ws.Range("K132").Value = 8771.571599999999
$ws.Range("L132").Value = 8821.636200000001
$ws.Range("M132").Value = -6241.571599999999
$ws.Range("N132").Value = -13881.6362

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("H32").Value = 606.5
$ws.Range("I32").Value = 606.5
$ws.Range("K32").Value = 606.5
$ws.Range("M32").Value = -289.5
$ws.Range("N26").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5874.9165
$ws.Range("I62").Value = 6620
$ws.Range("J62").Value = 5626.5557
$ws.Range("K62").Value = 6620
$ws.Range("L62").Value = 5626.5557
$ws.Range("M62").Value = -5996
$ws.Range("N62").Value = -6874.5557
$ws.Range("H65").Value = 5874.9165
$ws.Range("I65").Value = 6620
$ws.Range("J65").Value = 5626.5557
$ws.Range("K65").Value = 33100
$ws.Range("L65").Value = 28132.7785
$ws.Range("M65").Value = -29980
$ws.Range("N65").Value = -34372.7785
$ws.Range("H107").Value = 933.3333
$ws.Range("I107").Value = 800
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 2400
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -480
$ws.Range("N107").Value = -6840
$ws.Range("H141").Value = 36623.89
$ws.Range("J141").Value = 36623.89
$ws.Range("L141").Value = 36623.89
$ws.Range("N141").Value = -46983.89
